# Update the Gantt planner worksheet:
#  - change the scrolled/selected view (topLeftCell + selection)
#  - update several "PLAN DURATION" (column D) values
#  - fill in rows 16 & 18 (columns C/D/E/F/G) which were previously zeroed out

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("PLAN DURATION") value changes
$ws.Range("D5").Value  = 3
$ws.Range("D6").Value  = 4
$ws.Range("D7").Value  = 3
$ws.Range("D8").Value  = 4
$ws.Range("D9").Value  = 3
$ws.Range("D10").Value = 5
$ws.Range("D11").Value = 5
$ws.Range("D12").Value = 2
$ws.Range("D13").Value = 3
$ws.Range("D14").Value = 5
$ws.Range("D17").Value = 4
$ws.Range("D19").Value = 2

# Row 16 - fill in plan/actual data
$ws.Range("C16").Value = 8
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = 7
$ws.Range("F16").Value = 5
$ws.Range("G16").Value = 0.1

# Row 18 - fill in plan/actual data
$ws.Range("C18").Value = 10
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 11
$ws.Range("F18").Value = 3
$ws.Range("G18").Value = 1

# Update the view: scroll position (top-left visible cell -> B6) and the
# active selection (-> E17)
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("E17").Select()
